$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3745.158
$ws.Range("I64").Value = 3781.32
$ws.Range("J64").Value = 3675.6155
$ws.Range("K64").Value = 3781.32
$ws.Range("L64").Value = 3675.6155
$ws.Range("M64").Value = -3533.32
$ws.Range("N64").Value = -4171.6155
$ws.Range("H67").Value = 3745.158
$ws.Range("I67").Value = 3781.32
$ws.Range("J67").Value = 3675.6155
$ws.Range("K67").Value = 3781.32
$ws.Range("L67").Value = 3675.6155
$ws.Range("M67").Value = -2923.32
$ws.Range("N67").Value = -5391.6155
$ws.Range("H74").Value = 4605.1
$ws.Range("I74").Value = 4066.9167
$ws.Range("J74").Value = 5412.375
$ws.Range("K74").Value = 4066.9167
$ws.Range("L74").Value = 5412.375
$ws.Range("M74").Value = -3130.9167
$ws.Range("N74").Value = -7284.375
$ws.Range("H76").Value = 10527.6
$ws.Range("I76").Value = 19617.666
$ws.Range("K76").Value = 19617.666
$ws.Range("M76").Value = -19302.666
$ws.Range("H77").Value = 4605.1
$ws.Range("I77").Value = 4066.9167
$ws.Range("J77").Value = 5412.375
$ws.Range("K77").Value = 20334.5835
$ws.Range("L77").Value = 27061.875
$ws.Range("M77").Value = -15654.5835
$ws.Range("N77").Value = -36421.875
$ws.Range("H79").Value = 10527.6
$ws.Range("I79").Value = 19617.666
$ws.Range("K79").Value = 19617.666
$ws.Range("M79").Value = -18525.666
$ws.Range("H98").Value = 1738.6923
$ws.Range("I98").Value = 1259.762
$ws.Range("K98").Value = 1259.762
$ws.Range("M98").Value = 238.2380000000001
$ws.Range("H122").Value = 1738.6923
$ws.Range("I122").Value = 1259.762
$ws.Range("K122").Value = 3779.286
$ws.Range("M122").Value = -1329.286
$ws.Range("H138").Value = 5165.3726
$ws.Range("I138").Value = 970.8214
$ws.Range("J138").Value = 10271.782
$ws.Range("K138").Value = 2912.4642
$ws.Range("L138").Value = 30815.346
$ws.Range("M138").Value = 2227.5358
$ws.Range("N138").Value = -41095.346
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5024.8667
$ws.Range("I32").Value = 3871
$ws.Range("K32").Value = 3871
$ws.Range("M32").Value = -3584
$ws.Range("H61").Value = 4891.3
$ws.Range("I61").Value = 4891.3
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4891.3
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4679.3
$ws.Range("N61").ClearContents()
$ws.Range("H63").Value = 83335720
$ws.Range("I63").Value = 111113650
$ws.Range("J63").Value = 1950
$ws.Range("K63").Value = 111113650
$ws.Range("L63").Value = 1950
$ws.Range("M63").Value = -111112964
$ws.Range("N63").Value = -3322
$ws.Range("H66").Value = 83335720
$ws.Range("I66").Value = 111113650
$ws.Range("J66").Value = 1950
$ws.Range("K66").Value = 555568250
$ws.Range("L66").Value = 9750
$ws.Range("M66").Value = -555564818
$ws.Range("N66").Value = -16614
$ws.Range("H136").Value = 4891.3
$ws.Range("I136").Value = 4891.3
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 14673.9
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12123.9
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19609760
$ws.Range("I86").Value = 27779510
$ws.Range("K86").Value = 27779510
$ws.Range("M86").Value = -27778387
$ws.Range("H89").Value = 19609760
$ws.Range("I89").Value = 27779510
$ws.Range("K89").Value = 138897550
$ws.Range("M89").Value = -138891934
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5014.5713
$ws.Range("I62").Value = 5300.5
$ws.Range("J62").Value = 4299.75
$ws.Range("K62").Value = 5300.5
$ws.Range("L62").Value = 4299.75
$ws.Range("M62").Value = -4676.5
$ws.Range("N62").Value = -5547.75
$ws.Range("H65").Value = 5014.5713
$ws.Range("I65").Value = 5300.5
$ws.Range("J65").Value = 4299.75
$ws.Range("K65").Value = 26502.5
$ws.Range("L65").Value = 21498.75
$ws.Range("M65").Value = -23382.5
$ws.Range("N65").Value = -27738.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3460.9375
$ws.Range("J132").Value = 3491.6667
$ws.Range("L132").Value = 31425.0003
$ws.Range("N132").Value = -36485.0003
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6278.75
$ws.Range("I70").Value = 6294.5
$ws.Range("J70").Value = 6200
$ws.Range("K70").Value = 6294.5
$ws.Range("L70").Value = 6200
$ws.Range("M70").Value = -6024.5
$ws.Range("N70").Value = -6740
$ws.Range("H73").Value = 6278.75
$ws.Range("I73").Value = 6294.5
$ws.Range("J73").Value = 6200
$ws.Range("K73").Value = 6294.5
$ws.Range("L73").Value = 6200
$ws.Range("M73").Value = -5358.5
$ws.Range("N73").Value = -8072
$ws.Range("H80").Value = 2900
$ws.Range("I80").Value = 2750
$ws.Range("K80").Value = 2750
$ws.Range("M80").Value = -1752
$ws.Range("H83").Value = 2900
$ws.Range("I83").Value = 2750
$ws.Range("K83").Value = 13750
$ws.Range("M83").Value = -8758
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 39070.11
$ws.Range("I7").Value = 57160.723
$ws.Range("K7").Value = 57160.723
$ws.Range("M7").Value = -57048.723
$ws.Range("H40").Value = 100003530
$ws.Range("I40").Value = 100003530
$ws.Range("K40").Value = 100003530
$ws.Range("M40").Value = -100003394
$ws.Range("H69").Value = 98999.5
$ws.Range("J69").Value = 97999
$ws.Range("L69").Value = 97999
$ws.Range("N69").Value = -99621
$ws.Range("H72").Value = 98999.5
$ws.Range("J72").Value = 97999
$ws.Range("L72").Value = 293997
$ws.Range("N72").Value = -302109
$ws.Range("H81").Value = 42200
$ws.Range("J81").Value = 42200
$ws.Range("L81").Value = 42200
$ws.Range("N81").Value = -44196
$ws.Range("H82").Value = 80236.92
$ws.Range("I82").Value = 1330
$ws.Range("K82").Value = 1330
$ws.Range("M82").Value = -969
$ws.Range("H84").Value = 42200
$ws.Range("J84").Value = 42200
$ws.Range("L84").Value = 126600
$ws.Range("N84").Value = -136584
$ws.Range("H85").Value = 80236.92
$ws.Range("I85").Value = 1330
$ws.Range("K85").Value = 1330
$ws.Range("M85").Value = -82
$ws.Range("H126").Value = 39070.11
$ws.Range("I126").Value = 57160.723
$ws.Range("K126").Value = 171482.169
$ws.Range("M126").Value = -169012.169
